# Applies the "target function and limits" edit:
#  - inserts a "H2 Export Limit" column after the H2 columns
#  - inserts a "NH3 Export Limit" column after the NH3 columns
#  - appends a "CH3OH Export Limit" column at the end
# and fills in the corresponding values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before current column E (H2 Export Limit),
# shifting the NH3/CH3OH columns one to the right.
$ws.Range("E1").EntireColumn.Insert() | Out-Null

# Insert a new column before current column I (NH3 Export Limit),
# shifting the CH3OH columns one to the right.
$ws.Range("I1").EntireColumn.Insert() | Out-Null

# Headers for the newly inserted / appended columns
$ws.Range("E1").Value = "H2 Export Limit"
$ws.Range("I1").Value = "NH3 Export Limit"
$ws.Range("M1").Value = "CH3OH Export Limit"

# H2 Export Limit values (column E, rows 2-5)
$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 4
$ws.Range("E4").Value = 5
$ws.Range("E5").Value = 0

# NH3 Export Limit values (column I, rows 2-5)
$ws.Range("I2").Value = 5
$ws.Range("I3").Value = 6
$ws.Range("I4").Value = 2
$ws.Range("I5").Value = 4

# CH3OH Export Limit values (column M, rows 2-5)
$ws.Range("M2").Value = 8
$ws.Range("M3").Value = 1
$ws.Range("M4").Value = 3
$ws.Range("M5").Value = 4

# Size the new columns the same way Excel's "best fit" would after typing
# in the new header/values (matches the widths recorded in the workbook,
# rounded to this host's column-width pixel grid).
$ws.Columns.Item(5).ColumnWidth = 13.5
$ws.Columns.Item(9).ColumnWidth = 14.8333333333333
$ws.Columns.Item(13).ColumnWidth = 17.5

# Match the selection state recorded in the saved workbook
$ws.Range("F13").Select() | Out-Null
